# Update the "98_2" legislative-activity summary sheet:
#  - Pages of proceedings (House), C8: 12291 -> 12293
#  - Private bills enacted into law (Total), D11: 41 -> 46

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("98_2")

$ws.Range("C8").Value = 12293
$ws.Range("D11").Value = 46
